$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"-0.0002036580944324816"
$ws.Range("B2").Value = [double]"-3.03860379158778"
$ws.Range("C2").Value = [double]"-0.158562862399611"
$ws.Range("D2").Value = [double]"-0.2295784818484649"
$ws.Range("E2").Value = [double]"-8.828092784836333"
$ws.Range("F2").Value = [double]"-0.01102871272359567"
$ws.Range("G2").Value = [double]"-4.700623400836381E-06"
$ws.Range("H2").Value = [double]"-4.725445040508107E-06"
$ws.Range("I2").Value = [double]"-0.007043862153231701"
$ws.Range("J2").Value = [double]"-1.913075139828851"
$ws.Range("K2").Value = [double]"-2.863318232608441E-05"
$ws.Range("L2").Value = [double]"0.000692648549182994"
$ws.Range("M2").Value = [double]"-0.0002833098994483963"
$ws.Range("N2").Value = [double]"-2.358692474438168E-07"
$ws.Range("O2").Value = [double]"-7.002033470975084E-05"
$ws.Range("P2").Value = [double]"-0.0004926923049738186"
$ws.Range("Q2").Value = [double]"-0.0004978085136901718"
$ws.Range("R2").Value = [double]"-0.001299525934825169"
$ws.Range("S2").Value = [double]"-8.847708585665492E-09"
$ws.Range("T2").Value = [double]"-3.327072149998979E-06"
$ws.Range("U2").Value = [double]"-0.003800327442201475"
